# Update the cryptos list (Price/Volume columns) with refreshed quotes.
# Note: many Price values look numeric (e.g. "238.27") but must stay plain
# text, exactly like the source data, so a leading apostrophe is used to
# force a text value for those cells and avoid Excel auto-converting them
# to numbers (which would lose trailing zeros / use scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.270.23"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "1.885.54"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("D5").Value = "'238.27"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4673"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Value = "'0.06584"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").Value = "'19.89"
$ws.Range("E10").Value = "  +5.71%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "'98.04"
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "1.890.69"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "'5.121"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").Value = "'0.6692"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "'283.71"
$ws.Range("E16").Value = "  +11.16%  "
$ws.Range("D17").Value = "30.282.34"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "2.136.88"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "'12.61"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "'5.366"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "'0.000007305"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'6.175"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "'9.353"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'164.96"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'19.18"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "'1.993"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "'1.375"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'0.09735"
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("E31").Value = "  -5.21%  "
$ws.Range("D32").Value = "'1.483"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Value = "'4.180"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").Value = "'0.04696"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "'0.7102"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'2.709"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'0.01873"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "'6.655"
$ws.Range("E39").Value = "  +6.83%  "
$ws.Range("D40").Value = "'2.525"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").Value = "'72.49"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.974"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8698"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").Value = "'104.12"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").Value = "'1.0000"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "'985.09"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "'7.235"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").Value = "'9.241"
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").Value = "'34.10"
$ws.Range("E51").Value = "  -2.01%  "
